$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 0
$ws.Range("A5").Value = 3
